$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-08 Tuesday" "2025-07-09 Wednesday"

Replace-Text "359×3=1077" "935×5=4675"
Replace-Text "239×7=1673" "157×8=1256"
Replace-Text "647×2=1294" "769×7=5383"
Replace-Text "716×6=4296" "177×6=1062"
Replace-Text "299×3=897" "611×2=1222"

Replace-Text "338×3=1014" "498×9=4482"
Replace-Text "808×2=1616" "458×6=2748"
Replace-Text "731×8=5848" "892×8=7136"
Replace-Text "757×3=2271" "554×8=4432"
Replace-Text "962×6=5772" "222×9=1998"

Replace-Text "813×8=6504" "852×5=4260"
Replace-Text "401×4=1604" "343×3=1029"
Replace-Text "822×4=3288" "773×4=3092"
Replace-Text "239×6=1434" "699×3=2097"
Replace-Text "720×8=5760" "289×6=1734"

Replace-Text "833×9=7497" "603×8=4824"
Replace-Text "762×2=1524" "727×3=2181"
Replace-Text "450×7=3150" "952×3=2856"
Replace-Text "473×7=3311" "636×2=1272"
Replace-Text "852×4=3408" "414×8=3312"

Replace-Text "480×7=3360" "627×6=3762"
Replace-Text "766×3=2298" "384×3=1152"
Replace-Text "733×7=5131" "892×2=1784"
Replace-Text "586×2=1172" "602×6=3612"
Replace-Text "885×6=5310" "288×2=576"
